$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "Wins"/"Losses"/"Ties" columns, reusing the
# same formatting (bold, centered, bordered) as the existing header row by
# copying the format from the neighboring header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 79
    $ws.Cells.Item($row, 31).Value = 83
    $ws.Cells.Item($row, 32).Value = 0
}
